$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021年" / 132872 data row, matching the style used by the
# preceding year rows (bold, centered, bordered) by copying the formatting
# from the row above rather than constructing a brand-new style.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2021年"

$ws.Range("B12").Value = 132872
